# "added AIM 429 templates"
# Update the Channel Number values (column C) for the four data rows and
# leave the selection on C5, matching the authored template edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 26
$ws.Range("C3").Value = 26
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

$ws.Range("C5").Select()
